$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 15
$newRow = 16

# Copy the formatting (style) of column A from the row above, so the new
# row's A cell keeps the same border/alignment/bold style (s="1").
$ws.Cells.Item($srcRow, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = 14
$ws.Cells.Item($newRow, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item($newRow, 3).Value  = 1.071875655893269
$ws.Cells.Item($newRow, 4).Value  = 0.8773520961605883
$ws.Cells.Item($newRow, 5).Value  = 1.003181088641286
$ws.Cells.Item($newRow, 6).Value  = 1.071875655893269
$ws.Cells.Item($newRow, 7).Value  = 0.9206958191208094
$ws.Cells.Item($newRow, 8).Value  = 1.068054474609536
$ws.Cells.Item($newRow, 9).Value  = 1.02487355900856
$ws.Cells.Item($newRow, 10).Value = 0.8773520961605883
$ws.Cells.Item($newRow, 11).Value = 0.9402665924009372
$ws.Cells.Item($newRow, 12).Value = 1.006071124147103
$ws.Cells.Item($newRow, 13).Value = 0.9943387822390081
